$d = $word.ActiveDocument
$brk = [string][char]11

# 1. Update the date in the first paragraph
$d.Content.Find.Execute("June 17, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "June 21, 2022", 2) | Out-Null

# 2. Rework the second paragraph (liner description)
$p2 = $d.Paragraphs(2)
$p2.Range.Text = 'One(1) liner fabricated from ENTER MATERIAL NAME HERE'
$r = $p2.Range
$r.Collapse(0)
$r.InsertAfter($brk + $brk + '110''-0.0" diameter X 50''-0.0" deep')
$r.Collapse(0)
$r.InsertAfter(' with ENTER DEPTH EXTENSIONS HERE. ')

# 3. Insert the square-footage summary paragraph
$r = $p2.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.ParagraphFormat.LineSpacingRule = 2
$p3.Range.ParagraphFormat.LineSpacing = 24
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertAfter($brk + 'Bottom square footage:                                                                               12,432''')
$r3.Collapse(0)
$r3.InsertAfter($brk + 'Sidewall square footage:                                                                             ')
$r3.Collapse(0)
$r3.InsertAfter('17,602''')
$r3.Collapse(0)
$r3.InsertAfter($brk + 'Square footage:                                                                                              30,034''')
$r3.Collapse(0)
$r3.InsertAfter($brk + '5%:                                                                                                                     ')
$r3.Collapse(0)
$r3.InsertAfter('1,502''')
$r3.Font.Underline = 1
$r3.Collapse(0)
$r3.InsertAfter($brk + 'Total square footage:                                                                                   31,536''')
$r3.Font.Underline = 1
$r3.Collapse(0)
$r3.InsertAfter($brk + 'Cost of material:                                                                                              $2.00')

# 4. Insert the final cost paragraph
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$p4.Range.ParagraphFormat.LineSpacingRule = 2
$p4.Range.ParagraphFormat.LineSpacing = 24
$r4 = $p4.Range
$r4.Collapse(0)
$r4.InsertAfter($brk + 'Liner cost:                                                                                                    $63,071.40')
$r4.Collapse(0)
$r4.InsertAfter($brk + 'Total cost for one (1) lining system:                                                  $63,071.40')
$r4.Font.Underline = 1

Write-Output "done"
